$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 180; this shifts the existing rows 180-193
# down to 181-194 (and their formatting/styles travel with them).
$ws.Rows("180:180").Insert()

# Populate the newly inserted row 180 with the new weekly price record.
$ws.Range("A180").Value = 10
$ws.Range("B180").Value = "Vega Modelo de Temuco"
$ws.Range("C180").Value = "La Araucanía"
$ws.Range("D180").Value = 44826
$ws.Range("E180").Value = 9
$ws.Range("F180").Value = 100114007
$ws.Range("G180").Value = "Jengibre"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 200
$ws.Range("K180").Value = 20000
$ws.Range("L180").Value = 20000
$ws.Range("M180").Value = 20000
$ws.Range("N180").Value = '$/caja 13 kilos'
$ws.Range("O180").Value = "Perú"
$ws.Range("P180").Value = 1538
$ws.Range("Q180").Value = 13
$ws.Range("R180").Value = "Hortaliza"
